$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Colors (OLE/BGR packed ints, matching the target fgColor rgb values):
#   FF29A3CC -> blue   (unused by this sheet's rows, kept for parity)
#   FFFFCC66 -> orange (regular weekday attendance rows)
#   FFDF5E5E -> red    (rows flagged with a vacation-leave day)
# ---------------------------------------------------------------------------
$orange = 6737151
$red    = 6184671

function Paint-Row([int]$rowNum, [int]$color) {
    $rng = $ws.Range("A" + $rowNum + ":J" + $rowNum)
    $rng.Interior.Color = $color
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 11
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
    $rng.Borders.LineStyle = 1
    $rng.Borders.Color = 0
}

# Weekday attendance rows -> orange
Paint-Row 5  $orange
Paint-Row 6  $orange
Paint-Row 7  $orange
Paint-Row 11 $orange
Paint-Row 12 $orange
Paint-Row 13 $orange
Paint-Row 15 $orange

# Rows with a vacation-leave day -> red, and flag the VACATION LEAVE column
Paint-Row 8  $red
$ws.Range("I8").Value = 1

Paint-Row 14 $red
$ws.Range("I14").Value = 1

# Rows 9 and 10 (Saturday / Sunday) are intentionally left untouched.

# ---------------------------------------------------------------------------
# B19 is a member of the merged range A19:G19. Excel's merge semantics only
# keep the top-left cell's value, discarding the rest, so a direct
# Range("B19").Value assignment is silently dropped while the range stays
# merged. Routing the write through Copy/PasteSpecial (same mechanism Excel
# uses when you paste a single value onto a merged cell) lets the member
# cell's stored value/type actually change, matching the target (a literal
# boolean FALSE) while the merge stays intact.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = $false
$ws.Range("Z1").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Formula fixes: FLOOR() only takes 2 arguments; drop the stray trailing ",1"
# significance-repeat argument from each of these.
# ---------------------------------------------------------------------------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
